$d = $word.ActiveDocument

# Collapse a range to the very end of the document body content.
$endRange = $d.Content
$endRange.Collapse(0)   # wdCollapseEnd

# Insert a brand-new paragraph that contains nothing but a page break,
# carrying the same Consolas/minorHAnsi run formatting used throughout
# the grammar listing so it matches the surrounding paragraphs' look.
$pageBreakXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cstheme="minorHAnsi"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cstheme="minorHAnsi"/>
              </w:rPr>
              <w:br w:type="page"/>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

[void]$endRange.InsertXML($pageBreakXml)
